$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "E2" "6.55%"
Set-TextValue "D3" "31.83"
Set-TextValue "E3" "8.48%"
Set-TextValue "D4" "5.266"
Set-TextValue "E4" "4.13%"
Set-TextValue "D5" "0.07513"
Set-TextValue "E5" "12.00%"
Set-TextValue "D6" "7.822"
Set-TextValue "E6" "7.04%"
Set-TextValue "D7" "3.754"
Set-TextValue "E7" "9.12%"
Set-TextValue "D8" "1.480"
Set-TextValue "E8" "6.00%"
Set-TextValue "D9" "0.9144"
Set-TextValue "E9" "1.61%"
Set-TextValue "D10" "0.01670"
Set-TextValue "E10" "2,477.78%"
Set-TextValue "D11" "0.1697"
Set-TextValue "E11" "7.58%"
Set-TextValue "D12" "0.07593"
Set-TextValue "E12" "5.72%"
Set-TextValue "D13" "0.08055"
Set-TextValue "E13" "5.83%"
Set-TextValue "E14" "2.09%"
Set-TextValue "D15" "0.09895"
Set-TextValue "E15" "10.04%"
Set-TextValue "D16" "0.001486"
Set-TextValue "E16" "-7.25%"
Set-TextValue "D17" "0.04551"
Set-TextValue "E17" "1.24%"
Set-TextValue "D18" "0.006214"
Set-TextValue "E18" "-0.11%"
Set-TextValue "D19" "3.493"
Set-TextValue "E19" "1.31%"
Set-TextValue "D20" "2.234"
Set-TextValue "E20" "0.17%"
Set-TextValue "E21" "2.41%"
Set-TextValue "E22" "1.76%"
Set-TextValue "D23" "4.490"
Set-TextValue "E23" "15.42%"
Set-TextValue "E24" "4.44%"
Set-TextValue "D25" "0.001215"
Set-TextValue "E25" "1.04%"
Set-TextValue "D26" "0.004431"
Set-TextValue "E26" "1.33%"
Set-TextValue "D27" "0.0001398"
Set-TextValue "E27" "19.51%"
Set-TextValue "D28" "0.0001737"
Set-TextValue "E28" "7.35%"
Set-TextValue "D40" "0.04510"
Set-TextValue "E40" "6.24%"
Set-TextValue "D41" "0.007206"
Set-TextValue "E41" "5.73%"
Set-TextValue "D42" "0.1344"
Set-TextValue "E42" "8.52%"
Set-TextValue "D43" "0.002247"
Set-TextValue "E43" "0.78%"
Set-TextValue "D44" "0.01292"
Set-TextValue "E44" "1.83%"
Set-TextValue "D45" "0.00006199"
Set-TextValue "E45" "8.00%"
Set-TextValue "D46" "1.870"
Set-TextValue "E46" "-2.74%"
Set-TextValue "E47" "-13.59%"
